$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.324023666666666
$ws.Range("H2").Value = 3.972071
$ws.Range("I2").Value = 0.01518042398701374
$ws.Range("J2").Value = 0.01518042398701374
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.378475333333333
$ws.Range("N2").Value = 4.135426
$ws.Range("O2").Value = 0.05609715574531157
$ws.Range("P2").Value = 0.05609715574531156
$ws.Range("Q2").Value = 1.825133965249555
$ws.Range("R2").Value = 16.426205687246
$ws.Range("S2").Value = 0.0008515786086793733
$ws.Range("T2").Value = 0.0008515786086793732

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.324023666666666
$ws.Range("H3").Value = 3.972071
$ws.Range("I3").Value = 0.01518042398701374
$ws.Range("J3").Value = 0.01518042398701374
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.16176133333333
$ws.Range("N3").Value = 39.485284
$ws.Range("O3").Value = 0.5356188518899525
$ws.Range("P3").Value = 0.5356188518899525
$ws.Range("Q3").Value = 17.42648350035155
$ws.Range("R3").Value = 156.838351503164
$ws.Range("S3").Value = 0.008130921267126993
$ws.Range("T3").Value = 0.008130921267126993

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.324023666666666
$ws.Range("H4").Value = 3.972071
$ws.Range("I4").Value = 0.01518042398701374
$ws.Range("J4").Value = 0.01518042398701374
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.829094
$ws.Range("N4").Value = 2.487282
$ws.Range("O4").Value = 0.03374004171190829
$ws.Range("P4").Value = 0.03374004171190828
$ws.Range("Q4").Value = 1.097740077891333
$ws.Range("R4").Value = 9.879660701021999
$ws.Range("S4").Value = 0.0005121881385262967
$ws.Range("T4").Value = 0.0005121881385262966

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.324023666666666
$ws.Range("H5").Value = 3.972071
$ws.Range("I5").Value = 0.01518042398701374
$ws.Range("J5").Value = 0.01518042398701374
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.203668
$ws.Range("N5").Value = 27.611004
$ws.Range("O5").Value = 0.3745439506528278
$ws.Range("P5").Value = 0.3745439506528276
$ws.Range("Q5").Value = 12.18587425214267
$ws.Range("R5").Value = 109.672868269284
$ws.Range("S5").Value = 0.005685735972681076
$ws.Range("T5").Value = 0.005685735972681074

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 81.17653533333333
$ws.Range("H6").Value = 243.529606
$ws.Range("I6").Value = 0.9307191821270077
$ws.Range("J6").Value = 0.9307191821270075
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.378475333333333
$ws.Range("N6").Value = 4.135426
$ws.Range("O6").Value = 0.05609715574531157
$ws.Range("P6").Value = 0.05609715574531156
$ws.Range("Q6").Value = 111.8998516024618
$ws.Range("R6").Value = 1007.098664422156
$ws.Range("S6").Value = 0.05221069891492775
$ws.Range("T6").Value = 0.05221069891492774

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 81.17653533333333
$ws.Range("H7").Value = 243.529606
$ws.Range("I7").Value = 0.9307191821270077
$ws.Range("J7").Value = 0.9307191821270075
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.16176133333333
$ws.Range("N7").Value = 39.485284
$ws.Range("O7").Value = 0.5356188518899525
$ws.Range("P7").Value = 0.5356188518899525
$ws.Range("Q7").Value = 1068.426183924234
$ws.Range("R7").Value = 9615.835655318104
$ws.Range("S7").Value = 0.4985107397628235
$ws.Range("T7").Value = 0.4985107397628234

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.17653533333333
$ws.Range("H8").Value = 243.529606
$ws.Range("I8").Value = 0.9307191821270077
$ws.Range("J8").Value = 0.9307191821270075
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.829094
$ws.Range("N8").Value = 2.487282
$ws.Range("O8").Value = 0.03374004171190829
$ws.Range("P8").Value = 0.03374004171190828
$ws.Range("Q8").Value = 67.30297838565467
$ws.Range("R8").Value = 605.726805470892
$ws.Range("S8").Value = 0.0314025040270384
$ws.Range("T8").Value = 0.0314025040270384

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.17653533333333
$ws.Range("H9").Value = 243.529606
$ws.Range("I9").Value = 0.9307191821270077
$ws.Range("J9").Value = 0.9307191821270075
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.203668
$ws.Range("N9").Value = 27.611004
$ws.Range("O9").Value = 0.3745439506528278
$ws.Range("P9").Value = 0.3745439506528276
$ws.Range("Q9").Value = 747.1218805982694
$ws.Range("R9").Value = 6724.096925384424
$ws.Range("S9").Value = 0.3485952394222182
$ws.Range("T9").Value = 0.348595239422218

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.192675
$ws.Range("H10").Value = 3.578025
$ws.Range("I10").Value = 0.0136744626508778
$ws.Range("J10").Value = 0.0136744626508778
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.378475333333333
$ws.Range("N10").Value = 4.135426
$ws.Range("O10").Value = 0.05609715574531157
$ws.Range("P10").Value = 0.05609715574531156
$ws.Range("Q10").Value = 1.644073068183334
$ws.Range("R10").Value = 14.79665761365
$ws.Range("S10").Value = 0.0007670984610597382
$ws.Range("T10").Value = 0.0007670984610597379

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.192675
$ws.Range("H11").Value = 3.578025
$ws.Range("I11").Value = 0.0136744626508778
$ws.Range("J11").Value = 0.0136744626508778
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 13.16176133333333
$ws.Range("N11").Value = 39.485284
$ws.Range("O11").Value = 0.5356188518899525
$ws.Range("P11").Value = 0.5356188518899525
$ws.Range("Q11").Value = 15.69770369823333
$ws.Range("R11").Value = 141.2793332841
$ws.Range("S11").Value = 0.007324299985275205
$ws.Range("T11").Value = 0.007324299985275205

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.192675
$ws.Range("H12").Value = 3.578025
$ws.Range("I12").Value = 0.0136744626508778
$ws.Range("J12").Value = 0.0136744626508778
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.829094
$ws.Range("N12").Value = 2.487282
$ws.Range("O12").Value = 0.03374004171190829
$ws.Range("P12").Value = 0.03374004171190828
$ws.Range("Q12").Value = 0.9888396864500001
$ws.Range("R12").Value = 8.899557178050001
$ws.Range("S12").Value = 0.000461376940228549
$ws.Range("T12").Value = 0.0004613769402285489

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.192675
$ws.Range("H13").Value = 3.578025
$ws.Range("I13").Value = 0.0136744626508778
$ws.Range("J13").Value = 0.0136744626508778
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.203668
$ws.Range("N13").Value = 27.611004
$ws.Range("O13").Value = 0.3745439506528278
$ws.Range("P13").Value = 0.3745439506528276
$ws.Range("Q13").Value = 10.9769847319
$ws.Range("R13").Value = 98.7928625871
$ws.Range("S13").Value = 0.005121687264314312
$ws.Range("T13").Value = 0.00512168726431431

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.525915333333334
$ws.Range("H14").Value = 10.577746
$ws.Range("I14").Value = 0.04042593123510095
$ws.Range("J14").Value = 0.04042593123510094
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.378475333333333
$ws.Range("N14").Value = 4.135426
$ws.Range("O14").Value = 0.05609715574531157
$ws.Range("P14").Value = 0.05609715574531156
$ws.Range("Q14").Value = 4.860387314421778
$ws.Range("R14").Value = 43.743485829796
$ws.Range("S14").Value = 0.002267779760644714
$ws.Range("T14").Value = 0.002267779760644713

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.525915333333334
$ws.Range("H15").Value = 10.577746
$ws.Range("I15").Value = 0.04042593123510095
$ws.Range("J15").Value = 0.04042593123510094
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 13.16176133333333
$ws.Range("N15").Value = 39.485284
$ws.Range("O15").Value = 0.5356188518899525
$ws.Range("P15").Value = 0.5356188518899525
$ws.Range("Q15").Value = 46.40725609887378
$ws.Range("R15").Value = 417.665304889864
$ws.Range("S15").Value = 0.02165289087472694
$ws.Range("T15").Value = 0.02165289087472694

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.525915333333334
$ws.Range("H16").Value = 10.577746
$ws.Range("I16").Value = 0.04042593123510095
$ws.Range("J16").Value = 0.04042593123510094
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.829094
$ws.Range("N16").Value = 2.487282
$ws.Range("O16").Value = 0.03374004171190829
$ws.Range("P16").Value = 0.03374004171190828
$ws.Range("Q16").Value = 2.923315247374667
$ws.Range("R16").Value = 26.309837226372
$ws.Range("S16").Value = 0.001363972606115042
$ws.Range("T16").Value = 0.001363972606115042

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.525915333333334
$ws.Range("H17").Value = 10.577746
$ws.Range("I17").Value = 0.04042593123510095
$ws.Range("J17").Value = 0.04042593123510094
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.203668
$ws.Range("N17").Value = 27.611004
$ws.Range("O17").Value = 0.3745439506528278
$ws.Range("P17").Value = 0.3745439506528276
$ws.Range("Q17").Value = 32.45135412410934
$ws.Range("R17").Value = 292.062187116984
$ws.Range("S17").Value = 0.01514128799361426
$ws.Range("T17").Value = 0.01514128799361425
